# Update the "Rol" (role) values for each person listed in the
# "Seleccion Personal" worksheet. The write order below matches the
# order new shared-string entries must be appended in the workbook's
# string table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Seleccion Personal")

$ws.Range("C21").Value = "Desarrollador Analista Diseñador "
$ws.Range("F15").Value = "Desarrollador Analista Pruebas"
$ws.Range("F9").Value  = "Desarrollador Analista  Pruebas "
$ws.Range("C9").Value  = "Desarrollador  Analista  Diseñador"
$ws.Range("C15").Value = "Desarrollador   Diseñador  Analista "

# Reflect the saved view state (active cell/selection) recorded for this sheet.
$null = $ws.Range("L14").Select()
